$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between rows 116 and 117 ---
# (Indice/pais/torneio/temporada/data_partida in A:E stay on their own row)
$row116 = $ws.Range("F116:V116").Value2
$row117 = $ws.Range("F117:V117").Value2
$ws.Range("F116:V116").Value = $row117
$ws.Range("F117:V117").Value = $row116

# --- Rotate match data (columns F:V) among rows 146-149 ---
# new146 = old148, new147 = old149, new148 = old146, new149 = old147
$row146 = $ws.Range("F146:V146").Value2
$row147 = $ws.Range("F147:V147").Value2
$row148 = $ws.Range("F148:V148").Value2
$row149 = $ws.Range("F149:V149").Value2

$ws.Range("F146:V146").Value = $row148
$ws.Range("F147:V147").Value = $row149
$ws.Range("F148:V148").Value = $row146
$ws.Range("F149:V149").Value = $row147

# --- Append new row 151 (match Chindia Targoviste vs Concordia) ---
# Copy formatting from row 150 first so styles (bold/border index column,
# date number format column, etc.) match the rest of the table.
$ws.Range("A150:V150").Copy($ws.Range("A151:V151"))

$ws.Range("A151").Value = 150
$ws.Range("B151").Value = "romania"
$ws.Range("C151").Value = "liga-2"
$ws.Range("D151").Value = "2023-2024"
$ws.Range("E151").Value = 45263.4375
$ws.Range("F151").Value = "Chindia Targoviste"
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = "Concordia"
$ws.Range("I151").Value = 1
$ws.Range("J151").Value = 2.72
$ws.Range("K151").Value = "01/12/2023 22:42"
$ws.Range("L151").Value = 2.61
$ws.Range("M151").Value = "03/12/2023 10:15"
$ws.Range("N151").Value = 2.76
$ws.Range("O151").Value = "01/12/2023 22:42"
$ws.Range("P151").Value = 2.89
$ws.Range("Q151").Value = "03/12/2023 09:52"
$ws.Range("R151").Value = 2.76
$ws.Range("S151").Value = "01/12/2023 22:42"
$ws.Range("T151").Value = 3.05
$ws.Range("U151").Value = "03/12/2023 10:15"
$ws.Range("V151").Value = "https://www.betexplorer.com/football/romania/liga-2/chindia-targoviste-concordia/4vWhz2d9/"
